# The edit permutes the data rows 2-7 (each row's species/observation record)
# into a new order while leaving the columns that are identical across all
# six rows (C, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY, ...)
# untouched. Only columns A, B, D, E, F, G, H, I, M, Q, R actually differ
# between rows, so only those are read and rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "I", "M", "Q", "R")

# Snapshot the current ("before") values for rows 2..7 for the columns that
# vary, keyed by row number, before any writes happen.
$before = @{}
foreach ($r in 2..7) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowVals
}

# old row -> new row mapping (a permutation of rows 2..7)
$mapping = @{ 4 = 2; 7 = 3; 3 = 4; 2 = 5; 5 = 6; 6 = 7 }

foreach ($oldRow in $mapping.Keys) {
    $newRow = $mapping[$oldRow]
    $src = $before[$oldRow]
    foreach ($c in $cols) {
        $val = $src[$c]
        if ($val -eq $null) {
            $val = ""
        }
        $ws.Range("$c$newRow").Value = $val
    }
}
